$d = $word.ActiveDocument

# Replace the Jinja expression that calls name.full(middle="full")
# with the new name_full() call.
$d.Content.Find.Execute(
    "name.full(middle=""full"")",
    $true,   # MatchCase
    $false,  # MatchWholeWord
    $false,  # MatchWildcards
    $false,  # MatchSoundsLike
    $false,  # MatchAllWordForms
    $true,   # Forward
    1,       # Wrap (wdFindContinue)
    $false,  # Format
    "name_full()",
    2        # Replace (wdReplaceAll)
)
